$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# D-column values are prefixed with a leading apostrophe so Excel stores them
# as literal text instead of auto-converting number-like strings
# (e.g. "1.002", "0.07611", "22.372.22") into numeric values.

$ws.Range("D2").Value = "'22.372.22"
$ws.Range("E2").Value = "  -0.19%  "
$ws.Range("D3").Value = "'1.567.71"
$ws.Range("E3").Value = "  -0.29%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'1.002"
$ws.Range("E5").Value = "  +0.16%  "
$ws.Range("D6").Value = "'290.97"
$ws.Range("E6").Value = "  +0.41%  "
$ws.Range("D7").Value = "'0.3783"
$ws.Range("E7").Value = "  +3.04%  "
$ws.Range("D8").Value = "'49.17"
$ws.Range("E8").Value = "  -0.23%  "
$ws.Range("D9").Value = "'0.3405"
$ws.Range("E9").Value = "  +0.09%  "
$ws.Range("D10").Value = "'0.07611"
$ws.Range("D11").Value = "'1.139"
$ws.Range("E11").Value = "  -2.89%  "
$ws.Range("E12").Value = "  +0.23%  "
$ws.Range("D13").Value = "'21.06"
$ws.Range("E13").Value = "  -1.05%  "
$ws.Range("D14").Value = "'5.987"
$ws.Range("E14").Value = "  -1.35%  "
$ws.Range("D15").Value = "'6.928"
$ws.Range("E15").Value = "  +0.04%  "
$ws.Range("D16").Value = "'1.566.99"
$ws.Range("E16").Value = "  -0.24%  "
$ws.Range("D17").Value = "'0.00001134"
$ws.Range("E17").Value = "  -0.05%  "
$ws.Range("D18").Value = "'89.96"
$ws.Range("E18").Value = "  +0.02%  "
$ws.Range("D19").Value = "'0.06741"
$ws.Range("E19").Value = "  -0.05%  "
$ws.Range("E20").Value = "  +0.18%  "
$ws.Range("D21").Value = "'16.60"
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("D22").Value = "'6.196"
$ws.Range("E22").Value = "  -1.10%  "
$ws.Range("D23").Value = "'11.95"
$ws.Range("E23").Value = "  -0.51%  "
$ws.Range("D24").Value = "'22.354.15"
$ws.Range("E24").Value = "  -0.30%  "
$ws.Range("D25").Value = "'2.410"
$ws.Range("E25").Value = "  +1.02%  "
$ws.Range("D26").Value = "'2.696"
$ws.Range("E26").Value = "  -7.08%  "
$ws.Range("D27").Value = "'20.11"
$ws.Range("E27").Value = "  +0.26%  "
$ws.Range("D28").Value = "'147.36"
$ws.Range("E28").Value = "  +0.71%  "
$ws.Range("D29").Value = "'5.027"
$ws.Range("E29").Value = "  +0.87%  "
$ws.Range("D30").Value = "'126.06"
$ws.Range("E30").Value = "  +0.28%  "
$ws.Range("D31").Value = "'1.740.16"
$ws.Range("E31").Value = "  -0.34%  "
$ws.Range("D32").Value = "'2.016"
$ws.Range("E32").Value = "  -0.18%  "
$ws.Range("D33").Value = "'6.093"
$ws.Range("E33").Value = "  -2.49%  "
$ws.Range("D34").Value = "'0.9940"
$ws.Range("E34").Value = "  -3.44%  "
$ws.Range("D35").Value = "'10.11"
$ws.Range("E35").Value = "  -0.34%  "
$ws.Range("D36").Value = "'1.427"
$ws.Range("E36").Value = "  +9.39%  "
$ws.Range("D37").Value = "'0.08466"
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("D38").Value = "'0.02513"
$ws.Range("E38").Value = "  -1.11%  "
$ws.Range("D39").Value = "'0.2293"
$ws.Range("E39").Value = "  -1.33%  "
$ws.Range("D40").Value = "'0.06495"
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("D41").Value = "'5.410"
$ws.Range("E41").Value = "  -2.28%  "
$ws.Range("D42").Value = "'11.34"
$ws.Range("E42").Value = "  -3.43%  "
$ws.Range("D43").Value = "'0.6326"
$ws.Range("E43").Value = "  -0.62%  "
$ws.Range("E44").Value = "  +0.16%  "
$ws.Range("E45").Value = "  -2.01%  "
$ws.Range("D46").Value = "'3.808"
$ws.Range("E46").Value = "  +1.21%  "
$ws.Range("D47").Value = "'0.5937"
$ws.Range("E47").Value = "  -0.96%  "
$ws.Range("D48").Value = "'2.084"
$ws.Range("E48").Value = "  -1.34%  "
$ws.Range("D49").Value = "'1.255"
$ws.Range("E49").Value = "  -0.19%  "
$ws.Range("D50").Value = "'124.74"
$ws.Range("E50").Value = "  -0.32%  "
$ws.Range("D51").Value = "'0.07321"
$ws.Range("E51").Value = "  +0.38%  "
